$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 175.09091
$ws.Range("I39").Value = 81.40000000000001
$ws.Range("J39").Value = 253.16667
$ws.Range("K39").Value = 244.2
$ws.Range("L39").Value = 759.50001
$ws.Range("M39").Value = 51.79999999999998
$ws.Range("N39").Value = -1351.50001
$ws.Range("H40").Value = 2169.5
$ws.Range("I40").Value = 2191.5386
$ws.Range("J40").Value = 2128.5715
$ws.Range("K40").Value = 2191.5386
$ws.Range("L40").Value = 2128.5715
$ws.Range("M40").Value = -2016.5386
$ws.Range("N40").Value = -2478.5715
$ws.Range("H98").Value = 431453.97
$ws.Range("I98").Value = 621929.6
$ws.Range("J98").Value = 2883.75
$ws.Range("K98").Value = 621929.6
$ws.Range("L98").Value = 2883.75
$ws.Range("M98").Value = -620431.6
$ws.Range("N98").Value = -5879.75
$ws.Range("H116").Value = 11538106
$ws.Range("I116").Value = 17305782
$ws.Range("K116").Value = 17305782
$ws.Range("M116").Value = -17302340
$ws.Range("H122").Value = 431453.97
$ws.Range("I122").Value = 621929.6
$ws.Range("J122").Value = 2883.75
$ws.Range("K122").Value = 1865788.8
$ws.Range("L122").Value = 8651.25
$ws.Range("M122").Value = -1863338.8
$ws.Range("N122").Value = -13551.25
$ws.Range("H137").Value = 20834382
$ws.Range("I137").Value = 27778526
$ws.Range("J137").Value = 1955.0834
$ws.Range("K137").Value = 83335578
$ws.Range("L137").Value = 5865.2502
$ws.Range("M137").Value = -83333028
$ws.Range("N137").Value = -10965.2502

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 66061.875
$ws.Range("I2").Value = 115187.78
$ws.Range("K2").Value = 115187.78
$ws.Range("M2").Value = -115074.78
$ws.Range("H32").Value = 15782.421
$ws.Range("I32").Value = 2566.7385
$ws.Range("K32").Value = 2566.7385
$ws.Range("M32").Value = -2279.7385
$ws.Range("H46").Value = 7630.4
$ws.Range("J46").Value = 7630.4
$ws.Range("L46").Value = 7630.4
$ws.Range("N46").Value = -8268.4
$ws.Range("H61").Value = 1235.2094
$ws.Range("I61").Value = 1235.2094
$ws.Range("K61").Value = 1235.2094
$ws.Range("M61").Value = -1023.2094
$ws.Range("H74").Value = 4614
$ws.Range("I74").Value = 1546.6666
$ws.Range("J74").Value = 10529.571
$ws.Range("K74").Value = 1546.6666
$ws.Range("L74").Value = 10529.571
$ws.Range("M74").Value = -672.6666
$ws.Range("N74").Value = -12277.571
$ws.Range("H77").Value = 4614
$ws.Range("I77").Value = 1546.6666
$ws.Range("J77").Value = 10529.571
$ws.Range("K77").Value = 7733.333000000001
$ws.Range("L77").Value = 52647.855
$ws.Range("M77").Value = -3365.333000000001
$ws.Range("N77").Value = -61383.855
$ws.Range("H116").Value = 66061.875
$ws.Range("I116").Value = 115187.78
$ws.Range("K116").Value = 115187.78
$ws.Range("M116").Value = -112893.78
$ws.Range("H122").Value = 7680.35
$ws.Range("I122").Value = 8529.823
$ws.Range("J122").Value = 2866.6667
$ws.Range("K122").Value = 25589.469
$ws.Range("L122").Value = 8600.000100000001
$ws.Range("M122").Value = -23139.469
$ws.Range("N122").Value = -13500.0001
$ws.Range("H136").Value = 1235.2094
$ws.Range("I136").Value = 1235.2094
$ws.Range("K136").Value = 3705.6282
$ws.Range("M136").Value = -1155.6282

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 66061.875
$ws.Range("I3").Value = 115187.78
$ws.Range("K3").Value = 115187.78
$ws.Range("M3").Value = -115073.78
$ws.Range("H107").Value = 999.375
$ws.Range("I107").Value = 956.9231
$ws.Range("K107").Value = 956.9231
$ws.Range("M107").Value = 963.0769

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1072.375
$ws.Range("I16").Value = 479.8
$ws.Range("J16").Value = 2060
$ws.Range("K16").Value = 479.8
$ws.Range("L16").Value = 2060
$ws.Range("M16").Value = -192.8
$ws.Range("N16").Value = -2634
$ws.Range("H31").Value = 1915.4872
$ws.Range("I31").Value = 1120.9678
$ws.Range("K31").Value = 1120.9678
$ws.Range("M31").Value = -825.9677999999999
$ws.Range("H34").Value = 1915.4872
$ws.Range("I34").Value = 1120.9678
$ws.Range("K34").Value = 1120.9678
$ws.Range("M34").Value = -918.9677999999999
$ws.Range("H52").Value = 39700
$ws.Range("J52").Value = 39700
$ws.Range("L52").Value = 39700
$ws.Range("N52").Value = -40288
$ws.Range("H58").Value = 1791.2609
$ws.Range("I58").Value = 1029.3077
$ws.Range("J58").Value = 2781.8
$ws.Range("K58").Value = 1029.3077
$ws.Range("L58").Value = 2781.8
$ws.Range("M58").Value = -826.3077000000001
$ws.Range("N58").Value = -3187.8
$ws.Range("H113").Value = 1072.375
$ws.Range("I113").Value = 479.8
$ws.Range("J113").Value = 2060
$ws.Range("K113").Value = 479.8
$ws.Range("L113").Value = 2060
$ws.Range("M113").Value = 1690.2
$ws.Range("N113").Value = -6400
$ws.Range("H122").Value = 1572
$ws.Range("I122").Value = 1465
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4395
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -1945
$ws.Range("N122").Value = -10900
$ws.Range("H125").Value = 40042
$ws.Range("J125").Value = 40042
$ws.Range("L125").Value = 40042
$ws.Range("N125").Value = -44962
$ws.Range("H132").Value = 1564.6274
$ws.Range("I132").Value = 1419.6222
$ws.Range("K132").Value = 4258.8666
$ws.Range("M132").Value = -1728.8666
$ws.Range("H133").Value = 21134.857
$ws.Range("J133").Value = 22688
$ws.Range("L133").Value = 22688
$ws.Range("N133").Value = -27748
$ws.Range("H134").Value = 2050.7812
$ws.Range("I134").Value = 1344.1818
$ws.Range("K134").Value = 4032.5454
$ws.Range("M134").Value = -1497.5454
$ws.Range("H135").Value = 41499.168
$ws.Range("J135").Value = 41499.168
$ws.Range("L135").Value = 41499.168
$ws.Range("N135").Value = -51639.168
$ws.Range("H136").Value = 1791.2609
$ws.Range("I136").Value = 1029.3077
$ws.Range("J136").Value = 2781.8
$ws.Range("K136").Value = 3087.9231
$ws.Range("L136").Value = 8345.400000000001
$ws.Range("M136").Value = -537.9231
$ws.Range("N136").Value = -13445.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1303.3658
$ws.Range("I5").Value = 962.4167
$ws.Range("J5").Value = 1784.7059
$ws.Range("K5").Value = 2887.2501
$ws.Range("L5").Value = 5354.1177
$ws.Range("M5").Value = -2775.2501
$ws.Range("N5").Value = -5578.1177
$ws.Range("H122").Value = 650.93335
$ws.Range("J122").Value = 829.8
$ws.Range("L122").Value = 7468.2
$ws.Range("N122").Value = -12368.2
$ws.Range("H135").Value = 1303.3658
$ws.Range("I135").Value = 962.4167
$ws.Range("J135").Value = 1784.7059
$ws.Range("K135").Value = 8661.7503
$ws.Range("L135").Value = 16062.3531
$ws.Range("M135").Value = -6126.7503
$ws.Range("N135").Value = -21132.3531

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7196.45
$ws.Range("I70").Value = 7984.5454
$ws.Range("J70").Value = 6233.222
$ws.Range("K70").Value = 7984.5454
$ws.Range("L70").Value = 6233.222
$ws.Range("M70").Value = -7714.5454
$ws.Range("N70").Value = -6773.222
$ws.Range("H73").Value = 7196.45
$ws.Range("I73").Value = 7984.5454
$ws.Range("J73").Value = 6233.222
$ws.Range("K73").Value = 7984.5454
$ws.Range("L73").Value = 6233.222
$ws.Range("M73").Value = -7048.5454
$ws.Range("N73").Value = -8105.222
$ws.Range("H102").Value = 2879.7222
$ws.Range("I102").Value = 2766.818
$ws.Range("J102").Value = 3057.1428
$ws.Range("K102").Value = 2766.818
$ws.Range("L102").Value = 3057.1428
$ws.Range("M102").Value = -1144.818
$ws.Range("N102").Value = -6301.1428
$ws.Range("H113").Value = 1954.6428
$ws.Range("I113").Value = 1580.8125
$ws.Range("K113").Value = 1580.8125
$ws.Range("M113").Value = 589.1875
$ws.Range("H122").Value = 484178.47
$ws.Range("I122").Value = 556630.25
$ws.Range("J122").Value = 1166.6666
$ws.Range("K122").Value = 1669890.75
$ws.Range("L122").Value = 3499.9998
$ws.Range("M122").Value = -1667440.75
$ws.Range("N122").Value = -8399.9998
$ws.Range("H126").Value = 2467.3684
$ws.Range("I126").Value = 2195
$ws.Range("J126").Value = 2540
$ws.Range("K126").Value = 6585
$ws.Range("L126").Value = 7620
$ws.Range("M126").Value = -4115
$ws.Range("N126").Value = -12560
$ws.Range("H132").Value = 2927.2126
$ws.Range("I132").Value = 2596.3
$ws.Range("K132").Value = 7788.900000000001
$ws.Range("M132").Value = -5258.900000000001
$ws.Range("H134").Value = 15599.6
$ws.Range("J134").Value = 15599.6
$ws.Range("L134").Value = 46798.8
$ws.Range("N134").Value = -51868.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1908.0714
$ws.Range("I46").Value = 1594.5555
$ws.Range("J46").Value = 2472.4
$ws.Range("K46").Value = 1594.5555
$ws.Range("L46").Value = 2472.4
$ws.Range("M46").Value = -1406.5555
$ws.Range("N46").Value = -2848.4
$ws.Range("H122").Value = 3527.7778
$ws.Range("I122").Value = 2100
$ws.Range("J122").Value = 3611.7646
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 10835.2938
$ws.Range("M122").Value = -3850
$ws.Range("N122").Value = -15735.2938

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3269181.5
$ws.Range("J107").Value = 1276.6666
$ws.Range("L107").Value = 3829.9998
$ws.Range("N107").Value = -7669.9998
$ws.Range("H122").Value = 64552.625
$ws.Range("J122").Value = 2590
$ws.Range("L122").Value = 7770
$ws.Range("N122").Value = -12670
$ws.Range("H126").Value = 92119.82000000001
$ws.Range("I126").Value = 167685.5
$ws.Range("J126").Value = 1441
$ws.Range("K126").Value = 503056.5
$ws.Range("L126").Value = 4323
$ws.Range("M126").Value = -500586.5
$ws.Range("N126").Value = -9263
$ws.Range("H136").Value = 17245.209
$ws.Range("I136").Value = 18212.725
$ws.Range("J136").Value = 3216.25
$ws.Range("K136").Value = 54638.175
$ws.Range("L136").Value = 9648.75
$ws.Range("M136").Value = -52088.175
$ws.Range("N136").Value = -14748.75
